# The "cccd" value for the first student (row 2, column B) is updated from
# the plain number 111 to the text value "11123" (an 11-character CCCD
# number must be stored as text, not as a numeric value, since a numeric
# value would lose context / leading characters).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$cell = $ws.Range("B2")

# Mark the cell as Text first so Excel stores the assigned value as a
# string (t="s" in the XML) instead of coercing "11123" into a number.
$cell.NumberFormat = "@"
$cell.Value = "11123"

# Restore the cell's style back to the workbook default (same as before
# the edit) now that the value has already been recorded as text.
$cell.Style = "Normal"
